$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 8 through 16 (no longer present in the new data set)
$ws.Range("A8:D16").EntireRow.Delete() | Out-Null

# Update header row
$ws.Range("B1").Value = "Sizes"
$ws.Range("C1").Value = "My project"
$ws.Range("D1").Value = "enumerative_backtracking_solver.py"

# Copy style of the existing header cells (A1/B1) onto the new header cells
$ws.Range("A1:B1").Copy() | Out-Null
$ws.Range("C1:D1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Row 2 - test70419.txt
$ws.Range("A2").Value = "test70419.txt"
$ws.Range("B2").Value = "21*25"
$ws.Range("C2").Value = 3.737349033355713
$ws.Range("D2").Value = 64.78196930885315

# Row 3 - test70402.txt
$ws.Range("A3").Value = "test70402.txt"
$ws.Range("B3").Value = "24*18"
$ws.Range("C3").Value = 1.997181177139282
$ws.Range("D3").Value = 21.54502701759338

# Row 4 - test70399.txt
$ws.Range("A4").Value = "test70399.txt"
$ws.Range("B4").Value = "18*29"
$ws.Range("C4").Value = 1.694288015365601
$ws.Range("D4").Value = 336.943610906601

# Row 5 - test70446.txt
$ws.Range("A5").Value = "test70446.txt"
$ws.Range("B5").Value = "22*45"
$ws.Range("C5").Value = 2.374522924423218
$ws.Range("D5").Value = 141.6387090682983

# Row 6 - test70470.txt
$ws.Range("A6").Value = "test70470.txt"
$ws.Range("B6").Value = "45*41"
$ws.Range("C6").Value = 4.309165954589844
$ws.Range("D6").Value = 468.6386382579803

# Row 7 - test70468.txt
$ws.Range("A7").Value = "test70468.txt"
$ws.Range("B7").Value = "45*45"
$ws.Range("C7").Value = 6.089066982269287
$ws.Range("D7").Value = 321.1132562160492
